# Update the Vtn-Itgav L-R pair worksheet to include the "ECs" sending
# cluster (previously only a target cluster), per Dr Hou's advice.
# This rewrites the 9-row cartesian product of
# Sending cluster x Target cluster in {ECs, FAPs, sCs} for Ligand=Vtn,
# Receptor=Itgav, with refreshed expression/specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$rows = @(
    @("ECs", "Vtn", "Itgav", "ECs", 2, 0.6666666666666666, 7.134618, 21.403854, 0.0965317920926077, 0.0965317920926077, 3, 1, 16.535604, 49.606812, 0.2120453146491552, 0.2120453146491552, 117.975217939272, 1061.776961453448, 0.02046911422792383, 0.02046911422792384),
    @("ECs", "Vtn", "Itgav", "FAPs", 2, 0.6666666666666666, 7.134618, 21.403854, 0.0965317920926077, 0.0965317920926077, 3, 1, 40.62063066666667, 121.861892, 0.5209011059384622, 0.5209011059384622, 289.812682725752, 2608.314144531768, 0.05028351725926105, 0.05028351725926105),
    @("ECs", "Vtn", "Itgav", "sCs", 2, 0.6666666666666666, 7.134618, 21.403854, 0.0965317920926077, 0.0965317920926077, 3, 1, 20.825229, 62.475687, 0.2670535794123827, 0.2670535794123827, 148.580053677522, 1337.220483097698, 0.02577916060542283, 0.02577916060542283),
    @("FAPs", "Vtn", "Itgav", "ECs", 3, 1, 17.50798033333334, 52.52394100000001, 0.2368839813846793, 0.2368839813846794, 3, 1, 16.535604, 49.606812, 0.2120453146491552, 0.2120453146491552, 289.505029631788, 2605.545266686092, 0.05023013836805896, 0.05023013836805897),
    @("FAPs", "Vtn", "Itgav", "FAPs", 3, 1, 17.50798033333334, 52.52394100000001, 0.2368839813846793, 0.2368839813846794, 3, 1, 40.62063066666667, 121.861892, 0.5209011059384622, 0.5209011059384622, 711.1852028395971, 6400.666825556374, 0.1233931278823856, 0.1233931278823856),
    @("FAPs", "Vtn", "Itgav", "sCs", 3, 1, 17.50798033333334, 52.52394100000001, 0.2368839813846793, 0.2368839813846794, 3, 1, 20.825229, 62.475687, 0.2670535794123827, 0.2670535794123827, 364.607699769163, 3281.469297922467, 0.06326071513423484, 0.06326071513423485),
    @("sCs", "Vtn", "Itgav", "ECs", 3, 1, 49.26691733333334, 147.800752, 0.6665842265227129, 0.666584226522713, 3, 1, 16.535604, 49.606812, 0.2120453146491552, 0.2120453146491552, 814.6582353247361, 7331.924117922625, 0.1413460620531724, 0.1413460620531725),
    @("sCs", "Vtn", "Itgav", "FAPs", 3, 1, 49.26691733333334, 147.800752, 0.6665842265227129, 0.666584226522713, 3, 1, 40.62063066666667, 121.861892, 0.5209011059384622, 0.5209011059384622, 2001.253253082532, 18011.27927774279, 0.3472244607968155, 0.3472244607968156),
    @("sCs", "Vtn", "Itgav", "sCs", 3, 1, 49.26691733333334, 147.800752, 0.6665842265227129, 0.666584226522713, 3, 1, 20.825229, 62.475687, 0.2670535794123827, 0.2670535794123827, 1025.994835590736, 9233.953520316625, 0.178013703672725, 0.178013703672725)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $data[$i, $j] = $row[$j]
    }
}

$ws.Range("A2:T10").Value2 = $data
